$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new item ("FLAGYL 125MG/5ML 100 ML SUSPENSION") was inserted alphabetically
# between item 1 (BETADERM) and the former item 2 (TERRAMYCIN). That pushes the
# old "totals" row (was row 9) and the footer row (was row 10) down by one, and
# the former item-2 row's data (TERRAMYCIN) now becomes item 3 on the new row 9,
# while row 8 (item 2) is updated in place to the new FLAGYL data.

# Insert a new row at 9 - shifts old row 9 (totals) -> 10, old row 10 (footer) -> 11
$ws.Rows("9:9").Insert()

# ---- New row 9: item 3, carries what used to be item 2's data (TERRAMYCIN) ----
$ws.Range("A9").Value = 3
$ws.Range("C9").Value = "TERRAMYCIN EYE OINT. 5 GM"
$ws.Range("H9").Value = "6:0"
$ws.Range("L9").Value = "1"
$ws.Range("N9").Value = "28.00"
$ws.Range("P9").Value = "28.0000"
$ws.Range("Q9").Value = "1:0"

# Merge cells for row 9 (matches the pattern used by rows 7 and 8)
$ws.Range("A9:B9").Merge()
$ws.Range("C9:G9").Merge()
$ws.Range("H9:K9").Merge()
$ws.Range("L9:M9").Merge()
$ws.Range("N9:O9").Merge()
$ws.Rows("9:9").RowHeight = 25.5

# ---- Row 8 (item 2) now becomes the new FLAGYL entry ----
$ws.Range("C8").Value = "FLAGYL 125MG/5ML 100 ML SUSPENSION"
$ws.Range("H8").Value = "5:0"
$ws.Range("N8").Value = "26.00"
$ws.Range("P8").Value = "26.0000"

# ---- Totals row (shifted from row 9 to row 10): update transaction count ----
$ws.Range("N10").Value = 72
$ws.Rows("10:10").RowHeight = 25.5

Write-Output "done"
